# Post processing / Excel pre-process update
# - Add path to optimum (pareto) front json file used for post processing
# - Bump population size and number of generations
# - Active sheet moves from "Range Variables" back to "Project"

$wb = $excel.ActiveWorkbook

# --- "Project" sheet: add a new row with the path to the optimum front data ---
$wsProject = $wb.Worksheets.Item("Project")
$wsProject.Range("A8").Value = "path_opt_front "
$wsProject.Range("B8").Value = "C:\Users\jon\git\deap1\examples\ga\pareto_front\zdt1_front.json"

# --- "Parameters" sheet: larger population / more generations ---
$wsParameters = $wb.Worksheets.Item("Parameters")
$wsParameters.Range("B2").Value = 100
$wsParameters.Range("B3").Value = 1000
$wsParameters.Activate()
$wsParameters.Range("B4").Select()

# --- Make "Project" the active sheet again (selection moved to D13) ---
$wsProject.Activate()
$wsProject.Range("D13").Select()
